$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column B. This pushes the existing
# B:E (Jun_17 / Jun_15 / Jun_13 / Jun_10+ratings) columns to E:H,
# mirroring how the new "10th" snapshot columns were added in front of
# the existing weekly columns.
$ws.Columns("B:D").Insert()

# New header row for the three freshly inserted snapshot columns.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Every existing analyst row (2-27) gets "UN" placeholders in the three
# new columns, just like the rest of the table before any rating data
# has been collected for that date.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Two new analyst groups/firms added at the bottom of the table.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
